# Fruta / hortaliza, semanal
# The weekly refresh re-shuffled the already-existing daily observations
# (rows 2..25) across the date grid. Columns A,B,C,E,F,G,H,I,J are
# constants for this market/product and do not change; columns
# D,K,L,M,N,O,P,Q,R,S,T get redistributed among rows according to the
# mapping below (targetRow -> sourceRow, using the ORIGINAL/before values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move together as a single observation "record"
$cols = @("D","K","L","M","N","O","P","Q","R","S","T")

# target row -> source row (both refer to the workbook BEFORE this edit)
$rowMap = @{
    2  = 19
    3  = 20
    4  = 21
    5  = 4
    6  = 5
    7  = 24
    8  = 25
    9  = 11
    10 = 17
    11 = 18
    12 = 7
    13 = 8
    14 = 22
    15 = 23
    16 = 15
    17 = 16
    18 = 12
    19 = 13
    20 = 14
    21 = 9
    22 = 10
    23 = 6
    24 = 2
    25 = 3
}

# 1) Snapshot the current ("before") values for every row/column we need,
#    so that writes to one row never affect values later read from another.
#    (Value2 is used instead of Value because it reliably returns plain
#    numbers/strings/date-serials.)
$snapshot = @{}
foreach ($row in 2..25) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowData
}

# 2) Write the redistributed values into every target row.
foreach ($targetRow in 2..25) {
    $sourceRow = $rowMap[$targetRow]
    $sourceData = $snapshot[$sourceRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value2 = $sourceData[$col]
    }
}
